$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15; this shifts the existing rows 15-37 down to 16-38
# (matches the "Product Details Page" test case being inserted into the sheet).
$ws.Rows("15:15").Insert()

# The freshly inserted row doesn't carry the table's usual cell formatting
# (border/font/alignment), so copy the formatting from the row directly below
# (which held the content that used to be row 15) onto the new row 15.
$ws.Range("A16:F16").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Product Details Page" test case row.
$ws.Range("A15").Value = "Product Details Page "
$ws.Range("B15").Value = "TC3"
$ws.Range("C15").Value = "Fetching all product details "
$ws.Range("D15").Value = "Medium"
$ws.Range("E15").Value = "Sanity"

# Reflect the active cell/selection shown in the edited workbook.
$ws.Range("A15").Select() | Out-Null
